# repull data, push all data, mean calculation
# Update column F ("dSF") values on Sheet1 to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 4
    5  = 7
    7  = -1
    11 = 3
    21 = -1
    22 = 7
    23 = 5
    25 = -3
    26 = -6
    38 = -1
    43 = 1
    46 = -5
    51 = -3
    57 = 8
    60 = 8
    61 = -2
    62 = -2
    63 = 0
    68 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
